# Auto-generated script applying numeric corrections to market-data sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 1339.5385
$ws.Range("I9").Value = 291.4
$ws.Range("K9").Value = 291.4
$ws.Range("M9").Value = -122.4
# Row 96
$ws.Range("H96").Value = 1764
$ws.Range("I96").Value = 668.25
$ws.Range("J96").Value = 2859.75
$ws.Range("K96").Value = 2004.75
$ws.Range("L96").Value = 8579.25
$ws.Range("M96").Value = -631.75
$ws.Range("N96").Value = -11325.25
# Row 98
$ws.Range("H98").Value = 1719.8182
$ws.Range("I98").Value = 1719.8182
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1719.8182
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -221.8181999999999
$ws.Range("N98").ClearContents()
# Row 116
$ws.Range("H116").Value = 6200.6665
$ws.Range("J116").Value = 5766.6665
$ws.Range("L116").Value = 5766.6665
$ws.Range("N116").Value = -12650.6665
# Row 122
$ws.Range("H122").Value = 1719.8182
$ws.Range("I122").Value = 1719.8182
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5159.4546
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2709.4546
$ws.Range("N122").ClearContents()
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
# Row 132
$ws.Range("H132").Value = 1905.9
$ws.Range("I132").Value = 1820.0625
$ws.Range("K132").Value = 5460.1875
$ws.Range("M132").Value = -2930.1875
# Row 138
$ws.Range("H138").Value = 13182.147
$ws.Range("I138").Value = 11498
$ws.Range("J138").Value = 13233.182
$ws.Range("K138").Value = 34494
$ws.Range("L138").Value = 39699.546
$ws.Range("M138").Value = -29354
$ws.Range("N138").Value = -49979.546

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9108.806
$ws.Range("I32").Value = 8654.799999999999
$ws.Range("K32").Value = 8654.799999999999
$ws.Range("M32").Value = -8367.799999999999
# Row 61
$ws.Range("H61").Value = 2399.5386
$ws.Range("I61").Value = 1799.3334
$ws.Range("J61").Value = 3750
$ws.Range("K61").Value = 1799.3334
$ws.Range("L61").Value = 3750
$ws.Range("M61").Value = -1587.3334
$ws.Range("N61").Value = -4174
# Row 132
$ws.Range("H132").Value = 1677.2273
$ws.Range("I132").Value = 1279.3529
$ws.Range("J132").Value = 3030
$ws.Range("K132").Value = 3838.0587
$ws.Range("L132").Value = 9090
$ws.Range("M132").Value = -1308.0587
$ws.Range("N132").Value = -14150
# Row 136
$ws.Range("H136").Value = 2399.5386
$ws.Range("I136").Value = 1799.3334
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 5398.0002
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -2848.0002
$ws.Range("N136").Value = -16350

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1679.2
$ws.Range("I86").Value = 1632
$ws.Range("K86").Value = 1632
$ws.Range("M86").Value = -509
# Row 89
$ws.Range("H89").Value = 1679.2
$ws.Range("I89").Value = 1632
$ws.Range("K89").Value = 8160
$ws.Range("M89").Value = -2544
# Row 99
$ws.Range("H99").Value = 1635.6364
$ws.Range("I99").Value = 1499.375
$ws.Range("K99").Value = 1499.375
$ws.Range("M99").Value = -1.375

$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 4917.7144
$ws.Range("I94").Value = 5252.5
$ws.Range("K94").Value = 5252.5
$ws.Range("M94").Value = -4801.5
# Row 132
$ws.Range("H132").Value = 4454.4546
$ws.Range("I132").Value = 3755.4443
$ws.Range("K132").Value = 11266.3329
$ws.Range("M132").Value = -8736.332900000001
# Row 141
$ws.Range("H141").Value = 576655.7
$ws.Range("J141").Value = 576655.7
$ws.Range("L141").Value = 576655.7
$ws.Range("N141").Value = -587015.7

$ws = $wb.Worksheets.Item("CUL")
# Row 64
$ws.Range("H64").Value = 2000
$ws.Range("I64").Value = 2000
$ws.Range("K64").Value = 6000
$ws.Range("M64").Value = -5730
# Row 67
$ws.Range("H67").Value = 2000
$ws.Range("I67").Value = 2000
$ws.Range("K67").Value = 6000
$ws.Range("M67").Value = -5064
# Row 129
$ws.Range("H129").Value = 1792.9
$ws.Range("J129").Value = 3339.8
$ws.Range("L129").Value = 10019.4
$ws.Range("N129").Value = -20019.4

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 14289532
$ws.Range("J70").Value = 4503.5
$ws.Range("L70").Value = 4503.5
$ws.Range("N70").Value = -5043.5
# Row 73
$ws.Range("H73").Value = 14289532
$ws.Range("J73").Value = 4503.5
$ws.Range("L73").Value = 4503.5
$ws.Range("N73").Value = -6375.5
# Row 103
$ws.Range("H103").Value = 80000
$ws.Range("J103").Value = 80000
$ws.Range("L103").Value = 80000
$ws.Range("N103").Value = -82344
# Row 122
$ws.Range("H122").Value = 6438
$ws.Range("J122").Value = 4998.5
$ws.Range("L122").Value = 14995.5
$ws.Range("N122").Value = -19895.5
# Row 126
$ws.Range("H126").Value = 250003040
$ws.Range("I126").Value = 333334880
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 1000004640
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -1000002170
$ws.Range("N126").Value = -27440
# Row 132
$ws.Range("H132").Value = 2897.111
$ws.Range("I132").Value = 2449.8667
$ws.Range("J132").Value = 5133.3335
$ws.Range("K132").Value = 7349.6001
$ws.Range("L132").Value = 15400.0005
$ws.Range("M132").Value = -4819.6001
$ws.Range("N132").Value = -20460.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1136.75
$ws.Range("I22").Value = 757.1429000000001
$ws.Range("J22").Value = 1668.2
$ws.Range("K22").Value = 757.1429000000001
$ws.Range("L22").Value = 1668.2
$ws.Range("M22").Value = -462.1429000000001
$ws.Range("N22").Value = -2258.2
# Row 27
$ws.Range("H27").Value = 1136.75
$ws.Range("I27").Value = 757.1429000000001
$ws.Range("J27").Value = 1668.2
$ws.Range("K27").Value = 757.1429000000001
$ws.Range("L27").Value = 1668.2
$ws.Range("M27").Value = -650.1429000000001
$ws.Range("N27").Value = -1882.2
# Row 122
$ws.Range("H122").Value = 8427.444
$ws.Range("I122").Value = 8480.875
$ws.Range("K122").Value = 25442.625
$ws.Range("M122").Value = -22992.625
# Row 132
$ws.Range("H132").Value = 3876
$ws.Range("I132").Value = 1834.6666
$ws.Range("K132").Value = 5503.9998
$ws.Range("M132").Value = -2973.9998
# Row 136
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 406.25
$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 801.5
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 801.5
$ws.Range("M2").Value = 101
$ws.Range("N2").Value = -1025.5
# Row 4
$ws.Range("H4").Value = 3501
$ws.Range("I4").Value = 3501
$ws.Range("J4").Value = 3501
$ws.Range("K4").Value = 3501
$ws.Range("L4").Value = 3501
$ws.Range("M4").Value = -3388
$ws.Range("N4").Value = -3727
# Row 113
$ws.Range("H113").Value = 911.125
$ws.Range("I113").Value = 648
$ws.Range("J113").Value = 1490
$ws.Range("K113").Value = 1944
$ws.Range("L113").Value = 4470
$ws.Range("M113").Value = 226
$ws.Range("N113").Value = -8810
